# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) counts per game row, replacing the old "Strike#" derived values
$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 3
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 2
    15 = 1
    16 = 2
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
